# Apply the two changes described by the commit:
#  1. Slide 6's table switches to a different (built-in) table style.
#  2. The presentation's design/theme switches its colour palette from
#     the "Integral" theme to the default "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        # Table styles can't be assigned through the Style property
        # directly - PowerPoint requires ApplyStyle(styleId).
        $shp.Table.ApplyStyle("{6C4761D7-2364-41FC-81EF-4FD71553EBEA}")
    }
}

# --- 2. Swap the active theme's colour scheme (Integral -> Office) -------
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
